# Completed Automation:
# 1) Rename "LoginTest" sheet to "LoginTestSuccessFull" and fix the title cell / data
#    so every row uses the successful password.
# 2) Duplicate that sheet right after itself to create "LoginTestUnSuccessFull" that
#    holds the same InputUsername/InputPassword/ExpectedUsername table but without the
#    title row and with the "unsuccessful" password value in every data row.

$wb = $excel.ActiveWorkbook

# --- 1) Rename + fix up the original "LoginTest" sheet -----------------------------
$ws1 = $wb.Worksheets.Item("LoginTest")
$ws1.Name = "LoginTestSuccessFull"
$ws1.Range("A1").Value = "LoginTestSuccessFull"

# Make sure every data row carries the successful password (row 5 previously held the
# "unsuccessful" value by mistake).
$ws1.Range("B3").Value = "N0rthg4t31"
$ws1.Range("B4").Value = "N0rthg4t31"
$ws1.Range("B5").Value = "N0rthg4t31"

# --- 2) Duplicate it to build the "LoginTestUnSuccessFull" sheet -------------------
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "LoginTestUnSuccessFull"

# Drop the title row so the table starts at row 1 with the header.
$ws2.Rows.Item(1).Delete()

# Every data row uses the "unsuccessful" password this time.
$ws2.Range("B2").Value = "N0rthg4t311"
$ws2.Range("B3").Value = "N0rthg4t311"
$ws2.Range("B4").Value = "N0rthg4t311"

# Keep the first sheet as the active/selected one, matching the original workbook.
$ws1.Activate()
